$p = $ppt.ActivePresentation
Write-Host "HasTitleMaster: $($p.HasTitleMaster)"
Write-Host "HasNotesMaster: $($p.HasNotesMaster)"
Write-Host "HasHandoutMaster: $($p.HasHandoutMaster)"
try {
  $hm = $p.HandoutMaster
  Write-Host "HandoutMaster OK: $hm"
} catch {
  Write-Host "HandoutMaster error: $_"
}
